# infraAuto/silt-template.xlsx -- feat: api rest integration
# Renames the sheet and adds two new "armazem" / "integracao_via_servico_rest"
# columns (with their values) to the entity row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab (Entidade -> entidade)
$ws.Name = "entidade"

# New header cells (row 1) -- plain, unstyled like the existing E1 header
$ws.Range("P1").Value = "armazem"
$ws.Range("Q1").Value = "integracao_via_servico_rest"

# New data cells (row 2) -- centered like the rest of the data row
$dataRange = $ws.Range("P2:Q2")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$ws.Range("P2").Value = "IFC EMBU"
$ws.Range("Q2").Value = "sim"

# Match the new column widths as closely as this engine's width model allows
$ws.Columns.Item(16).ColumnWidth = 11
$ws.Columns.Item(17).ColumnWidth = 25.333333333333332

# The active cell ends up on the newly entered value
$ws.Range("Q1").Select()
